$d = $word.ActiveDocument

# Locate the "LOB1037: ..." requirement paragraph; the three paragraphs that
# follow it (a blank paragraph, the "Ver no Jupiter ..." line and the
# "© 2020 ..." footer line) are removed by this edit.
$count = $d.Paragraphs.Count
$anchor = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*LOB1037*") {
        $anchor = $i
        break
    }
}

if ($anchor -gt 0) {
    # Delete from the last of the three paragraphs back to the first so that
    # earlier indices stay valid while we work.
    $d.Paragraphs.Item($anchor + 3).Range.Delete()
    $d.Paragraphs.Item($anchor + 2).Range.Delete()
    $d.Paragraphs.Item($anchor + 1).Range.Delete()
}
